# Generate Report for Handback
# - Marks the 712ff180 and acc33caa source files as handed back (zh-cn and de-de)
# - Updates the Overview status + the per-language "Latest Target/Handback File"
#   and "Latest Handback DateTime" columns
# - Adds hyperlinks on the newly-populated "Latest Target File" cells
# - Widens the Status column(s) to fit the new, longer status text

$wb = $excel.ActiveWorkbook

$ovw = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$zhcnHandbackTime = "2017-02-28 08:03:07"
$dedeHandbackTime = "2017-02-28 08:03:28"

# ---------------------------------------------------------------------------
# Overview sheet: status columns for the first file (712ff180) flip to
# "Handed back: in sync with en-US" for both locales.
# ---------------------------------------------------------------------------
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 2 (712ff180) and 3 (acc33caa) get their target/handback
# file + datetime populated now that the handback round-tripped.
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("J2").Value = "712ff180-b74f-4ba8-95d1-44e62587c07f.md"
$zhcn.Range("K2").Value = "712ff180-b74f-4ba8-95d1-44e62587c07f.378470fdbe8e8c4e47be0e61c069a9614b96cd90.zh-cn.xlf"
$zhcn.Range("L2").Value = $zhcnHandbackTime

$zhcn.Range("J3").Value = "acc33caa-6a9d-47e3-a2a0-22020a6a2c37.md"
$zhcn.Range("K3").Value = "acc33caa-6a9d-47e3-a2a0-22020a6a2c37.c3e729d8a296e168435f8af3a10100dea1d787d5.zh-cn.xlf"
$zhcn.Range("L3").Value = $zhcnHandbackTime

$zhcn.Hyperlinks.Add($zhcn.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/1653045dadf229584bb191495158da850d74c077/e2e/712ff180-b74f-4ba8-95d1-44e62587c07f.md", "", "", "712ff180-b74f-4ba8-95d1-44e62587c07f.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/1653045dadf229584bb191495158da850d74c077/e2e/acc33caa-6a9d-47e3-a2a0-22020a6a2c37.md", "", "", "acc33caa-6a9d-47e3-a2a0-22020a6a2c37.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet: same shape of update, using the de-de xlf / handback datetime.
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("J2").Value = "712ff180-b74f-4ba8-95d1-44e62587c07f.md"
$dede.Range("K2").Value = "712ff180-b74f-4ba8-95d1-44e62587c07f.378470fdbe8e8c4e47be0e61c069a9614b96cd90.de-de.xlf"
$dede.Range("L2").Value = $dedeHandbackTime

$dede.Range("J3").Value = "acc33caa-6a9d-47e3-a2a0-22020a6a2c37.md"
$dede.Range("K3").Value = "acc33caa-6a9d-47e3-a2a0-22020a6a2c37.c3e729d8a296e168435f8af3a10100dea1d787d5.de-de.xlf"
$dede.Range("L3").Value = $dedeHandbackTime

$dede.Hyperlinks.Add($dede.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/b50341033f0fd15fa6274f1ba398f84eb0b4314c/e2e/712ff180-b74f-4ba8-95d1-44e62587c07f.md", "", "", "712ff180-b74f-4ba8-95d1-44e62587c07f.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/b50341033f0fd15fa6274f1ba398f84eb0b4314c/e2e/acc33caa-6a9d-47e3-a2a0-22020a6a2c37.md", "", "", "acc33caa-6a9d-47e3-a2a0-22020a6a2c37.md") | Out-Null

# ---------------------------------------------------------------------------
# Column widths widen to fit "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$ovw.Columns.Item(5).ColumnWidth = 29.1666666666667
$ovw.Columns.Item(6).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
